$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save" - copy G1's formatting (bold header style) then set value
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data column values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
